$wb = $excel.ActiveWorkbook

function Add-SensorRow($ws, $row, $timeVal, $idHex, $actualHex, $checksumHex, $totalDec, $idDec, $actualDec, $checksumDec) {
    $ws.Cells.Item($row, 1).Value2 = $timeVal
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = "0x01,0x90"
    $ws.Cells.Item($row, 3).Value = $idHex
    $ws.Cells.Item($row, 4).Value = $actualHex
    $ws.Cells.Item($row, 5).Value = $checksumHex

    $ws.Cells.Item($row, 6).Value2 = $totalDec
    $ws.Cells.Item($row, 7).Value2 = [double]$idDec
    $ws.Cells.Item($row, 8).Value2 = $actualDec
    $ws.Cells.Item($row, 9).Value2 = $checksumDec
}

$idDecBig = "5.686312626471138e+23"

$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
Add-SensorRow $ws1 29 45729.73239443287 "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0xd" 400 $idDecBig 400 13
Add-SensorRow $ws1 30 45729.73241640046 "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0xd" 400 $idDecBig 400 13
Add-SensorRow $ws1 31 45729.73243972223 "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0xd" 400 $idDecBig 400 13

$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
Add-SensorRow $ws2 29 45729.58037369213 "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x86," "0x4" 400 $idDecBig 390 4
Add-SensorRow $ws2 30 45729.58039555555 "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x86," "0x4" 400 $idDecBig 390 4
Add-SensorRow $ws2 31 45729.58041870371 "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x86," "0x4" 400 $idDecBig 390 4
